# Adds two new trailing columns, I ("I0") and J ("IF"), to the existing
# per-game stats table on Sheet1, extending the used range from A1:H53
# to A1:J53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Give the two new header cells the same look as the rest of the header
# row (bold font, thin border, centered/top aligned) by copying the
# formatting from the existing "IP" header cell (H1) before stamping in
# the new labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-53) ----------------------------------------------
# Each row is [row number, I0 value, IF value].
$rows = @(
    @(2,9,9),
    @(3,8,8),
    @(4,6,7),
    @(5,9,10),
    @(6,8,8),
    @(7,8,8),
    @(8,8,8),
    @(9,8,8),
    @(10,7,7),
    @(11,8,8),
    @(12,8,8),
    @(13,8,8),
    @(14,8,8),
    @(15,8,8),
    @(16,10,10),
    @(17,5,6),
    @(18,7,7),
    @(19,10,10),
    @(20,7,7),
    @(21,7,7),
    @(22,7,7),
    @(23,10,10),
    @(24,8,9),
    @(25,7,7),
    @(26,7,7),
    @(27,5,5),
    @(28,7,7),
    @(29,7,7),
    @(30,9,9),
    @(31,8,8),
    @(32,6,6),
    @(33,10,10),
    @(34,5,6),
    @(35,9,9),
    @(36,8,8),
    @(37,7,7),
    @(38,8,8),
    @(39,7,7),
    @(40,8,8),
    @(41,4,4),
    @(42,9,9),
    @(43,7,8),
    @(44,7,7),
    @(45,9,9),
    @(46,7,7),
    @(47,6,7),
    @(48,6,6),
    @(49,7,7),
    @(50,7,7),
    @(51,5,5),
    @(52,5,5),
    @(53,5,5)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 9).Value  = $r[1]   # column I -> I0
    $ws.Cells.Item($rowNum, 10).Value = $r[2]   # column J -> IF
}
